$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("posts")

# --- New column N: "IS_GENERATE" header + FALSE for every data row ---

# Header cell (N1): copy the format of the adjacent header cell (M1) so it
# picks up the existing bold/fill/border header style instead of minting a
# new cellXf, then set its text.
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)
$ws.Range("N1").Value = "IS_GENERATE"

# Data cells (N2:N12): copy the format of an existing boolean column cell
# (I2, style shared by the whole data block) and fill with FALSE.
$ws.Range("I2").Copy()
$ws.Range("N2:N12").PasteSpecial(-4122)
$ws.Range("N2:N12").Value = $false

# Clear the clipboard marching ants / clipboard reference.
$excel.CutCopyMode = 0

# --- Update the active selection to match the saved view state ---
$ws.Range("M12").Select() | Out-Null
